$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 319.77777
$ws.Range("J2").Value = 300
$ws.Range("L2").Value = 300
$ws.Range("N2").Value = -526

$ws.Range("H4").Value = 627
$ws.Range("I4").Value = 627
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 627
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -513
$ws.Range("N4").Value = ""

$ws.Range("H28").Value = 590.53845
$ws.Range("I28").Value = 495.9
$ws.Range("K28").Value = 495.9
$ws.Range("M28").Value = -10.89999999999998

$ws.Range("H33").Value = 240.92308
$ws.Range("I33").Value = 217.14285
$ws.Range("K33").Value = 217.14285
$ws.Range("M33").Value = 11.85714999999999

$ws.Range("H42").Value = 75
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 75
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 225
$ws.Range("M42").Value = ""
$ws.Range("N42").Value = -685

$ws.Range("H43").Value = 3639.1667
$ws.Range("I43").Value = 2944.6667
$ws.Range("J43").Value = 4333.6665
$ws.Range("K43").Value = 2944.6667
$ws.Range("L43").Value = 4333.6665
$ws.Range("M43").Value = -2875.6667
$ws.Range("N43").Value = -4471.6665

$ws.Range("H64").Value = 12380
$ws.Range("I64").Value = 4950
$ws.Range("J64").Value = 17333.334
$ws.Range("K64").Value = 4950
$ws.Range("L64").Value = 17333.334
$ws.Range("M64").Value = -4702
$ws.Range("N64").Value = -17829.334

$ws.Range("H67").Value = 12380
$ws.Range("I67").Value = 4950
$ws.Range("J67").Value = 17333.334
$ws.Range("K67").Value = 4950
$ws.Range("L67").Value = 17333.334
$ws.Range("M67").Value = -4092
$ws.Range("N67").Value = -19049.334

$ws.Range("H92").Value = 698.75
$ws.Range("I92").Value = 698.75
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 698.75
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 549.25
$ws.Range("N92").Value = ""

$ws.Range("H111").Value = 2887.25
$ws.Range("J111").Value = 2533
$ws.Range("L111").Value = 7599
$ws.Range("N111").Value = -13733

$ws.Range("H118").Value = 454.5
$ws.Range("I118").Value = 454.5
$ws.Range("K118").Value = 1363.5
$ws.Range("M118").Value = 293.5

$ws.Range("H132").Value = 6331.7144
$ws.Range("I132").Value = 3775
$ws.Range("K132").Value = 11325
$ws.Range("M132").Value = -8795

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 294.22223
$ws.Range("I22").Value = 318.5
$ws.Range("K22").Value = 318.5
$ws.Range("M22").Value = -145.5

$ws.Range("H86").Value = 2099.889
$ws.Range("I86").Value = 1985.7142
$ws.Range("J86").Value = 2499.5
$ws.Range("K86").Value = 1985.7142
$ws.Range("L86").Value = 2499.5
$ws.Range("M86").Value = -862.7141999999999
$ws.Range("N86").Value = -4745.5

$ws.Range("H89").Value = 2099.889
$ws.Range("I89").Value = 1985.7142
$ws.Range("J89").Value = 2499.5
$ws.Range("K89").Value = 9928.571
$ws.Range("L89").Value = 12497.5
$ws.Range("M89").Value = -4312.571
$ws.Range("N89").Value = -23729.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 942
$ws.Range("I16").Value = 930.4
$ws.Range("K16").Value = 930.4
$ws.Range("M16").Value = -643.4

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").Value = ""

$ws.Range("H63").Value = 99999
$ws.Range("J63").Value = 99999
$ws.Range("L63").Value = 99999
$ws.Range("N63").Value = -101371

$ws.Range("H66").Value = 99999
$ws.Range("J66").Value = 99999
$ws.Range("L66").Value = 299997
$ws.Range("N66").Value = -306861

$ws.Range("H95").Value = 35000
$ws.Range("J95").Value = 35000
$ws.Range("L95").Value = 35000
$ws.Range("N95").Value = -40492

$ws.Range("H105").Value = 1265.2632
$ws.Range("I105").Value = 844.9286
$ws.Range("K105").Value = 844.9286
$ws.Range("M105").Value = 902.0714

$ws.Range("H113").Value = 942
$ws.Range("I113").Value = 930.4
$ws.Range("K113").Value = 930.4
$ws.Range("M113").Value = 1239.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 600.7273
$ws.Range("I23").Value = 334.75
$ws.Range("J23").Value = 752.7143
$ws.Range("K23").Value = 1004.25
$ws.Range("L23").Value = 2258.1429
$ws.Range("M23").Value = -769.25
$ws.Range("N23").Value = -2728.1429

$ws.Range("H34").Value = 557.5
$ws.Range("I34").Value = 47.5
$ws.Range("J34").Value = 1067.5
$ws.Range("K34").Value = 142.5
$ws.Range("L34").Value = 3202.5
$ws.Range("M34").Value = -58.5
$ws.Range("N34").Value = -3370.5

$ws.Range("H68").Value = 1300
$ws.Range("I68").Value = 1300
$ws.Range("K68").Value = 3900
$ws.Range("M68").Value = -3089

$ws.Range("H71").Value = 1300
$ws.Range("I71").Value = 1300
$ws.Range("K71").Value = 11700
$ws.Range("M71").Value = -7644

$ws.Range("H86").Value = 1678
$ws.Range("I86").Value = 765
$ws.Range("J86").Value = 2069.2856
$ws.Range("K86").Value = 2295
$ws.Range("L86").Value = 6207.8568
$ws.Range("M86").Value = -1109
$ws.Range("N86").Value = -8579.856800000001

$ws.Range("H89").Value = 1678
$ws.Range("I89").Value = 765
$ws.Range("J89").Value = 2069.2856
$ws.Range("K89").Value = 6885
$ws.Range("L89").Value = 18623.5704
$ws.Range("M89").Value = -957
$ws.Range("N89").Value = -30479.5704

$ws.Range("H112").Value = 42049.1
$ws.Range("I112").Value = 1499
$ws.Range("J112").Value = 46554.668
$ws.Range("K112").Value = 4497
$ws.Range("L112").Value = 139664.004
$ws.Range("M112").Value = -3389
$ws.Range("N112").Value = -141880.004

$ws.Range("H113").Value = 231.25
$ws.Range("J113").Value = 254.6
$ws.Range("L113").Value = 763.8
$ws.Range("N113").Value = -5103.8

$ws.Range("H140").Value = 1948.1428
$ws.Range("J140").Value = 4987
$ws.Range("L140").Value = 14961
$ws.Range("N140").Value = -25321

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 826.38464
$ws.Range("I2").Value = 1131.4445
$ws.Range("J2").Value = 140
$ws.Range("K2").Value = 1131.4445
$ws.Range("L2").Value = 140
$ws.Range("M2").Value = -1018.4445
$ws.Range("N2").Value = -366

$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").Value = ""

$ws.Range("H57").Value = 21999.4
$ws.Range("I57").Value = 20000
$ws.Range("K57").Value = 20000
$ws.Range("M57").Value = -19180

$ws.Range("H107").Value = 955.6
$ws.Range("I107").Value = 659.3333
$ws.Range("K107").Value = 659.3333
$ws.Range("M107").Value = 1260.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2244.6924
$ws.Range("I22").Value = 2182.1667
$ws.Range("K22").Value = 2182.1667
$ws.Range("M22").Value = -1887.1667

$ws.Range("H27").Value = 2244.6924
$ws.Range("I27").Value = 2182.1667
$ws.Range("K27").Value = 2182.1667
$ws.Range("M27").Value = -2075.1667

$ws.Range("H68").Value = 2323.3333
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").Value = ""

$ws.Range("H71").Value = 2323.3333
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").Value = ""

$ws.Range("H140").Value = 79333.336
$ws.Range("J140").Value = 79333.336
$ws.Range("L140").Value = 79333.336
$ws.Range("N140").Value = -89693.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 50383
$ws.Range("J112").Value = 50383
$ws.Range("L112").Value = 50383
$ws.Range("N112").Value = -53337
